{"js": "// The document contains a paragraph whose visible text reads\n// \"<id>p022v_1</id>\" but whose content is split across five separate\n// runs: \"<id>\", \"p\", \"0\", \"22v_1\", \"</id>\" (leftover from an earlier\n// per-character edit). The commit collapses that back into a single\n// run carrying the text \"<id>p022v_1</id>\" (keeping the formatting of\n// the run that already held the \"<id>\"/\"</id>\" tags: Courier New,\n// color 7f6000, sz/szCs 18).\n//\n// Find the paragraph by its reconstructed text (stable whether or not\n// it is already split across runs) and rewrite its range in one shot;\n// Word's \"Replace\" insertion collapses the selection into a single run\n// using the formatting of the range's leading run, which is exactly\n// the target state.\n\nconst TARGET_TEXT = \"<id>p022v_1</id>\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === TARGET_TEXT) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(`Could not find paragraph with text ${JSON.stringify(TARGET_TEXT)}`);\n}\n\nconst range = target.getRange();\nrange.insertText(TARGET_TEXT, \"Replace\");\nawait context.sync();\n", "ps1": "# The document contains a paragraph whose visible text reads\n# \"<id>p022v_1</id>\" but whose content is split across five separate\n# runs: \"<id>\", \"p\", \"0\", \"22v_1\", \"</id>\" (leftover from an earlier\n# per-character edit). This collapses it back into a single run\n# carrying the text \"<id>p022v_1</id>\", taking on the formatting of\n# the run that already held the \"<id>\"/\"</id>\" tags (Courier New,\n# color 7f6000, sz/szCs 18) - exactly what Word's Find/Replace does\n# when the replacement text overwrites a multi-run match.\n#\n# \"<id>p022v_1</id>\" is unique across the whole document (the other\n# <id> tags wrap \"fig_p022v_1\" / \"fig_p022v_2\"), so a document-wide\n# Find & Replace unambiguously targets this one paragraph.\n\n$d = $word.ActiveDocument\n\n$needle = \"<id>p022v_1</id>\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = $needle\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = $needle\n\n# Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n#         MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n#         ReplaceWith, Replace)\n# Wrap:=1 (wdFindContinue), Replace:=2 (wdReplaceAll)\n$find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, $needle, 2)\n"}
